# Auto-generated Excel COM-interop script
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value2 = 37177
$ws.Range("J87").Value2 = 39569.332
$ws.Range("L87").Value2 = 39569.332
$ws.Range("N87").Value2 = -42065.332
$ws.Range("H90").Value2 = 37177
$ws.Range("J90").Value2 = 39569.332
$ws.Range("L90").Value2 = 118707.996
$ws.Range("N90").Value2 = -131187.996
$ws.Range("H138").Value2 = 2404.18
$ws.Range("I138").Value2 = 1129.4375
$ws.Range("J138").Value2 = 2646.988
$ws.Range("K138").Value2 = 3388.3125
$ws.Range("L138").Value2 = 7940.964
$ws.Range("M138").Value2 = 1751.6875
$ws.Range("N138").Value2 = -18220.964

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value2 = 27524
$ws.Range("I3").Value2 = 0
$ws.Range("J3").Value2 = 27524
$ws.Range("K3").Value2 = 0
$ws.Range("L3").Value2 = 27524
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value2 = -27754
$ws.Range("H6").Value2 = 34572.715
$ws.Range("I6").Value2 = 0
$ws.Range("J6").Value2 = 34572.715
$ws.Range("K6").Value2 = 0
$ws.Range("L6").Value2 = 34572.715
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value2 = -34918.715
$ws.Range("H8").Value2 = 70004.5
$ws.Range("I8").Value2 = 0
$ws.Range("J8").Value2 = 70004.5
$ws.Range("K8").Value2 = 0
$ws.Range("L8").Value2 = 70004.5
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value2 = -70292.5
$ws.Range("H11").Value2 = 2500490
$ws.Range("I11").Value2 = 5000000
$ws.Range("J11").Value2 = 980
$ws.Range("K11").Value2 = 5000000
$ws.Range("L11").Value2 = 980
$ws.Range("M11").Value2 = -4999856
$ws.Range("N11").Value2 = -1268

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value2 = 0
$ws.Range("I8").Value2 = 0
$ws.Range("K8").Value2 = 0
$ws.Range("M8").ClearContents()
$ws.Range("H70").Value2 = 75000
$ws.Range("J70").Value2 = 75000
$ws.Range("L70").Value2 = 75000
$ws.Range("N70").Value2 = -75586
$ws.Range("H73").Value2 = 75000
$ws.Range("J73").Value2 = 75000
$ws.Range("L73").Value2 = 75000
$ws.Range("N73").Value2 = -77028

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value2 = 1897000.2
$ws.Range("I6").Value2 = 3775250
$ws.Range("J6").Value2 = 18750.5
$ws.Range("K6").Value2 = 3775250
$ws.Range("L6").Value2 = 18750.5
$ws.Range("M6").Value2 = -3775137
$ws.Range("N6").Value2 = -18976.5
$ws.Range("H13").Value2 = 36001.668
$ws.Range("J13").Value2 = 36001.668
$ws.Range("L13").Value2 = 36001.668
$ws.Range("N13").Value2 = -36279.668
$ws.Range("H31").Value2 = 3813.3838
$ws.Range("I31").Value2 = 3211.5454
$ws.Range("J31").Value2 = 4443.881
$ws.Range("K31").Value2 = 3211.5454
$ws.Range("L31").Value2 = 4443.881
$ws.Range("M31").Value2 = -2916.5454
$ws.Range("N31").Value2 = -5033.881
$ws.Range("H34").Value2 = 3813.3838
$ws.Range("I34").Value2 = 3211.5454
$ws.Range("J34").Value2 = 4443.881
$ws.Range("K34").Value2 = 3211.5454
$ws.Range("L34").Value2 = 4443.881
$ws.Range("M34").Value2 = -3009.5454
$ws.Range("N34").Value2 = -4847.881
$ws.Range("H99").Value2 = 3637.3794
$ws.Range("I99").Value2 = 4063.818
$ws.Range("J99").Value2 = 3376.7778
$ws.Range("K99").Value2 = 4063.818
$ws.Range("L99").Value2 = 3376.7778
$ws.Range("M99").Value2 = -2565.818
$ws.Range("N99").Value2 = -6372.7778
$ws.Range("H122").Value2 = 512
$ws.Range("I122").Value2 = 512
$ws.Range("J122").Value2 = 0
$ws.Range("K122").Value2 = 1536
$ws.Range("L122").Value2 = 0
$ws.Range("M122").Value2 = 914
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value2 = 3637.3794
$ws.Range("I126").Value2 = 4063.818
$ws.Range("J126").Value2 = 3376.7778
$ws.Range("K126").Value2 = 12191.454
$ws.Range("L126").Value2 = 10130.3334
$ws.Range("M126").Value2 = -9721.454000000002
$ws.Range("N126").Value2 = -15070.3334
$ws.Range("H132").Value2 = 3437.6365
$ws.Range("I132").Value2 = 1720.4
$ws.Range("J132").Value2 = 4868.6665
$ws.Range("K132").Value2 = 5161.200000000001
$ws.Range("L132").Value2 = 14605.9995
$ws.Range("M132").Value2 = -2631.200000000001
$ws.Range("N132").Value2 = -19665.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value2 = 2961.9648
$ws.Range("I68").Value2 = 5635.05
$ws.Range("J68").Value2 = 1517.0541
$ws.Range("K68").Value2 = 16905.15
$ws.Range("L68").Value2 = 4551.1623
$ws.Range("M68").Value2 = -16094.15
$ws.Range("N68").Value2 = -6173.1623
$ws.Range("H71").Value2 = 2961.9648
$ws.Range("I71").Value2 = 5635.05
$ws.Range("J71").Value2 = 1517.0541
$ws.Range("K71").Value2 = 50715.45
$ws.Range("L71").Value2 = 13653.4869
$ws.Range("M71").Value2 = -46659.45
$ws.Range("N71").Value2 = -21765.4869

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value2 = 6250398.5
$ws.Range("I3").Value2 = 10000497
$ws.Range("J3").Value2 = 234.66667
$ws.Range("K3").Value2 = 10000497
$ws.Range("L3").Value2 = 234.66667
$ws.Range("M3").Value2 = -10000381
$ws.Range("N3").Value2 = -466.66667
$ws.Range("H7").Value2 = 8361.444
$ws.Range("I7").Value2 = 283.33334
$ws.Range("J7").Value2 = 12400.5
$ws.Range("K7").Value2 = 283.33334
$ws.Range("L7").Value2 = 12400.5
$ws.Range("M7").Value2 = -171.33334
$ws.Range("N7").Value2 = -12624.5
$ws.Range("H8").Value2 = 8361.444
$ws.Range("I8").Value2 = 283.33334
$ws.Range("J8").Value2 = 12400.5
$ws.Range("K8").Value2 = 283.33334
$ws.Range("L8").Value2 = 12400.5
$ws.Range("M8").Value2 = -144.33334
$ws.Range("N8").Value2 = -12678.5
$ws.Range("H11").Value2 = 9861107
$ws.Range("I11").Value2 = 15335445
$ws.Range("J11").Value2 = 7298
$ws.Range("K11").Value2 = 15335445
$ws.Range("L11").Value2 = 7298
$ws.Range("M11").Value2 = -15335306
$ws.Range("N11").Value2 = -7576
$ws.Range("H132").Value2 = 20356
$ws.Range("I132").Value2 = 1293.6571
$ws.Range("J132").Value2 = 53715.1
$ws.Range("K132").Value2 = 3880.9713
$ws.Range("L132").Value2 = 161145.3
$ws.Range("M132").Value2 = -1350.9713
$ws.Range("N132").Value2 = -166205.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value2 = 21500
$ws.Range("J92").Value2 = 21500
$ws.Range("L92").Value2 = 21500
$ws.Range("N92").Value2 = -26492
$ws.Range("H136").Value2 = 346374.72
$ws.Range("I136").Value2 = 589381.1
$ws.Range("J136").Value2 = 2115.6667
$ws.Range("K136").Value2 = 1768143.3
$ws.Range("L136").Value2 = 6347.000100000001
$ws.Range("M136").Value2 = -1765593.3
$ws.Range("N136").Value2 = -11447.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value2 = 1248.6666
$ws.Range("I4").Value2 = 646
$ws.Range("J4").Value2 = 1550
$ws.Range("K4").Value2 = 646
$ws.Range("L4").Value2 = 1550
$ws.Range("M4").Value2 = -533
$ws.Range("N4").Value2 = -1776
$ws.Range("H6").Value2 = 21113.666
$ws.Range("I6").Value2 = 1005
$ws.Range("J6").Value2 = 23627.25
$ws.Range("K6").Value2 = 1005
$ws.Range("L6").Value2 = 23627.25
$ws.Range("M6").Value2 = -890
$ws.Range("N6").Value2 = -23857.25
$ws.Range("H81").Value2 = 1628.76
$ws.Range("J81").Value2 = 1949.2307
$ws.Range("L81").Value2 = 3898.4614
$ws.Range("N81").Value2 = -6020.4614
$ws.Range("H84").Value2 = 1628.76
$ws.Range("J84").Value2 = 1949.2307
$ws.Range("L84").Value2 = 19492.307
$ws.Range("N84").Value2 = -30100.307
$ws.Range("H132").Value2 = 3713.5334
$ws.Range("I132").Value2 = 1016.90625
$ws.Range("J132").Value2 = 10351.385
$ws.Range("K132").Value2 = 3050.71875
$ws.Range("L132").Value2 = 31054.155
$ws.Range("M132").Value2 = -520.71875
$ws.Range("N132").Value2 = -36114.155
